# Updated via Streamlit Approval System
#
# Adds four new status columns (ACCEPTED / PAID / HOLD / REJECTED) to the
# pending-approval sheet, seeds boolean defaults for each data row, and
# backfills APPROVAL_1 / APPROVAL_2 ("HOLD") for the rows that were blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells AP1:AS1 -------------------------------------------
# Clone the formatting (bold, centered, bordered - style index 1) of the
# neighbouring header cell, then overwrite with the new header text.
$ws.Range("AO1").Copy($ws.Range("AP1"))
$ws.Range("AO1").Copy($ws.Range("AQ1"))
$ws.Range("AO1").Copy($ws.Range("AR1"))
$ws.Range("AO1").Copy($ws.Range("AS1"))

$ws.Range("AP1").Value = "ACCEPTED"
$ws.Range("AQ1").Value = "PAID"
$ws.Range("AR1").Value = "HOLD"
$ws.Range("AS1").Value = "REJECTED"

# --- Backfill APPROVAL_1 / APPROVAL_2 (columns AI / AJ) -----------------
# Rows that were blank become "HOLD".
$ws.Range("AI2:AJ2").Value = "HOLD"
$ws.Range("AI3:AJ3").Value = "HOLD"
$ws.Range("AI5:AJ5").Value = "HOLD"

# --- Seed the new boolean columns for each data row (2-16) --------------
$ws.Range("AP2:AQ4").Value = $false
$ws.Range("AR2:AR4").Value = $false
$ws.Range("AS2:AS4").Value = $false

$ws.Range("AP5:AQ5").Value = $false
$ws.Range("AR5").Value = $true
$ws.Range("AS5").Value = $false

$ws.Range("AP6:AQ16").Value = $false
$ws.Range("AR6:AR16").Value = $true
$ws.Range("AS6:AS16").Value = $false
